# qs-completingthesquare.docx refresh
#
# The original document had several short paragraphs whose text was
# split across many single-word runs (one <w:r> per word / space).
# This edit "flattens" each of those paragraphs back down to a single
# run containing the full sentence, without changing the visible text
# at all. We do this with Find & Replace across each paragraph's full
# text: Word's Find engine spans run boundaries when matching, and
# replacing the whole match collapses it into one run.

$d = $word.ActiveDocument

# Title: "Questions: Completing the square" (was split into 7 runs)
$d.Content.Find.Execute(
    "Questions: Completing the square", $true, $false, $false, $false,
    $false, $true, 1, $false, "Questions: Completing the square", 2) | Out-Null

# Author: "Tom Coleman" (was split into 3 runs)
$d.Content.Find.Execute(
    "Tom Coleman", $true, $false, $false, $false,
    $false, $true, 1, $false, "Tom Coleman", 2) | Out-Null

# Abstract: "A selection of questions for the study guide on completing
# the square." (was split into 23 runs, one per word/space)
$d.Content.Find.Execute(
    "A selection of questions for the study guide on completing the square.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on completing the square.",
    2) | Out-Null
